{"js": "// Office.js (Word JavaScript API) script.\n// Replaces the date title and all 100 arithmetic-problem cells in the\n// 10x5 table with their updated values, matching the target diff.\n// The mapping below is the ordered list of [oldText, newText] pairs,\n// one per paragraph in document order (paragraph 0 = the title line,\n// paragraphs 1..100 = the table cells in row-major order). Since every\n// paragraph in the body changes in this diff, we just walk body.paragraphs\n// in order and overwrite each one's text in place (preserving its run\n// formatting) rather than relying on text search, which keeps the edit\n// robust even though some \"old\" strings could coincidentally repeat.\n\nconst replacements = [[\"2023-04-10 Monday\", \"2023-04-11 Tuesday\"], [\"4+17=\", \"61+19=\"], [\"62-53=\", \"37+11=\"], [\"15+19=\", \"14+40=\"], [\"19+37=\", \"54+22=\"], [\"65-63=\", \"94-6=\"], [\"34+16=\", \"39+25=\"], [\"94-58=\", \"83-18=\"], [\"17+48=\", \"58-26=\"], [\"54+13=\", \"9+58=\"], [\"44-6=\", \"33+27=\"], [\"42+8=\", \"24+52=\"], [\"89-27=\", \"88-8=\"], [\"34-5=\", \"93-18=\"], [\"35-12=\", \"14-13=\"], [\"57-52=\", \"51-4=\"], [\"54-50=\", \"66+25=\"], [\"37+16=\", \"5+25=\"], [\"40-3=\", \"27+2=\"], [\"32+63=\", \"95-48=\"], [\"27+0=\", \"0+98=\"], [\"85+11=\", \"79+14=\"], [\"35+38=\", \"27-17=\"], [\"93-11=\", \"57-5=\"], [\"24+12=\", \"63+23=\"], [\"54-47=\", \"62+15=\"], [\"98-70=\", \"35+42=\"], [\"13+61=\", \"63+22=\"], [\"68-66=\", \"13+73=\"], [\"0+56=\", \"83-69=\"], [\"22+37=\", \"68-15=\"], [\"53-50=\", \"37-2=\"], [\"14+29=\", \"45+39=\"], [\"22+19=\", \"52-27=\"], [\"65-33=\", \"59-58=\"], [\"11+13=\", \"82+10=\"], [\"5+31=\", \"5+12=\"], [\"19+64=\", \"54+37=\"], [\"42+1=\", \"23+53=\"], [\"20+28=\", \"8+32=\"], [\"84-52=\", \"16+9=\"], [\"38+48=\", \"74+22=\"], [\"22-21=\", \"69+10=\"], [\"58-17=\", \"9+20=\"], [\"12+70=\", \"9+58=\"], [\"77-60=\", \"48-24=\"], [\"51+25=\", \"69+27=\"], [\"31+14=\", \"91-22=\"], [\"60+26=\", \"9+19=\"], [\"41+53=\", \"5+10=\"], [\"53+21=\", \"7+48=\"], [\"88-11=\", \"17+58=\"], [\"86-10=\", \"97-56=\"], [\"87-25=\", \"51-29=\"], [\"71-67=\", \"98-68=\"], [\"0+19=\", \"30+38=\"], [\"39+50=\", \"68+8=\"], [\"96-7=\", \"2+62=\"], [\"33+14=\", \"74-19=\"], [\"78-6=\", \"13+68=\"], [\"92-5=\", \"57-41=\"], [\"88-80=\", \"12+74=\"], [\"40-14=\", \"64+24=\"], [\"52+24=\", \"73-3=\"], [\"7-2=\", \"41-36=\"], [\"57-43=\", \"73-46=\"], [\"0+45=\", \"11-10=\"], [\"30-25=\", \"36-15=\"], [\"6+93=\", \"19+34=\"], [\"72-49=\", \"11+0=\"], [\"92-6=\", \"66-2=\"], [\"9+77=\", \"68-59=\"], [\"59+28=\", \"81-62=\"], [\"74-32=\", \"70-24=\"], [\"55+15=\", \"48-46=\"], [\"5+42=\", \"25+47=\"], [\"52-24=\", \"50-11=\"], [\"92-54=\", \"26+12=\"], [\"84-40=\", \"93-46=\"], [\"4+0=\", \"83-41=\"], [\"64+10=\", \"67-2=\"], [\"34-26=\", \"31-2=\"], [\"80-13=\", \"23+1=\"], [\"87-5=\", \"37+30=\"], [\"74-70=\", \"67+32=\"], [\"38-22=\", \"12+24=\"], [\"10+66=\", \"77-50=\"], [\"71+17=\", \"15+55=\"], [\"37-31=\", \"38+1=\"], [\"91-3=\", \"9+26=\"], [\"28-10=\", \"82-24=\"], [\"10+40=\", \"90-20=\"], [\"15-12=\", \"71-9=\"], [\"96-73=\", \"85-43=\"], [\"19-1=\", \"95-87=\"], [\"50-9=\", \"59-42=\"], [\"26-8=\", \"28+45=\"], [\"52+42=\", \"73+3=\"], [\"75-68=\", \"98-67=\"], [\"75+3=\", \"70-28=\"], [\"36-15=\", \"98-73=\"]];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = paragraphs.items[i];\n  para.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = paragraphs.items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Paragraph ${i} text mismatch: expected ${JSON.stringify(oldText)}, got ${JSON.stringify(para.text)}`\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Updates the date title paragraph and all 100 arithmetic-problem cells\n# in the 10x5 table (here laid out as a 20-row x 5-column Word table)\n# to match the target diff. Cells are addressed positionally via the\n# Tables/Rows/Columns object model (row-major order, same order as the\n# document/diff), and each cell's Range.Text is overwritten in place so\n# existing run formatting (font/size) on the paragraph is preserved.\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph (the date line) ---\n$titleOld = '2023-04-10 Monday'\n$titleNew = '2023-04-11 Tuesday'\n$titlePara = $d.Paragraphs.Item(1)\n$titleText = $titlePara.Range.Text.TrimEnd([char]13, [char]7)\nif ($titleText -ne $titleOld) {\n    throw \"Title paragraph text mismatch: expected [$titleOld], got [$titleText]\"\n}\n$titlePara.Range.Text = $titleNew\n\n# --- Table cells (row-major: row 1 col 1..5, row 2 col 1..5, ...) ---\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$cellValues = @(\n    @('4+17=', '61+19='),\n    @('62-53=', '37+11='),\n    @('15+19=', '14+40='),\n    @('19+37=', '54+22='),\n    @('65-63=', '94-6='),\n    @('34+16=', '39+25='),\n    @('94-58=', '83-18='),\n    @('17+48=', '58-26='),\n    @('54+13=', '9+58='),\n    @('44-6=', '33+27='),\n    @('42+8=', '24+52='),\n    @('89-27=', '88-8='),\n    @('34-5=', '93-18='),\n    @('35-12=', '14-13='),\n    @('57-52=', '51-4='),\n    @('54-50=', '66+25='),\n    @('37+16=', '5+25='),\n    @('40-3=', '27+2='),\n    @('32+63=', '95-48='),\n    @('27+0=', '0+98='),\n    @('85+11=', '79+14='),\n    @('35+38=', '27-17='),\n    @('93-11=', '57-5='),\n    @('24+12=', '63+23='),\n    @('54-47=', '62+15='),\n    @('98-70=', '35+42='),\n    @('13+61=', '63+22='),\n    @('68-66=', '13+73='),\n    @('0+56=', '83-69='),\n    @('22+37=', '68-15='),\n    @('53-50=', '37-2='),\n    @('14+29=', '45+39='),\n    @('22+19=', '52-27='),\n    @('65-33=', '59-58='),\n    @('11+13=', '82+10='),\n    @('5+31=', '5+12='),\n    @('19+64=', '54+37='),\n    @('42+1=', '23+53='),\n    @('20+28=', '8+32='),\n    @('84-52=', '16+9='),\n    @('38+48=', '74+22='),\n    @('22-21=', '69+10='),\n    @('58-17=', '9+20='),\n    @('12+70=', '9+58='),\n    @('77-60=', '48-24='),\n    @('51+25=', '69+27='),\n    @('31+14=', '91-22='),\n    @('60+26=', '9+19='),\n    @('41+53=', '5+10='),\n    @('53+21=', '7+48='),\n    @('88-11=', '17+58='),\n    @('86-10=', '97-56='),\n    @('87-25=', '51-29='),\n    @('71-67=', '98-68='),\n    @('0+19=', '30+38='),\n    @('39+50=', '68+8='),\n    @('96-7=', '2+62='),\n    @('33+14=', '74-19='),\n    @('78-6=', '13+68='),\n    @('92-5=', '57-41='),\n    @('88-80=', '12+74='),\n    @('40-14=', '64+24='),\n    @('52+24=', '73-3='),\n    @('7-2=', '41-36='),\n    @('57-43=', '73-46='),\n    @('0+45=', '11-10='),\n    @('30-25=', '36-15='),\n    @('6+93=', '19+34='),\n    @('72-49=', '11+0='),\n    @('92-6=', '66-2='),\n    @('9+77=', '68-59='),\n    @('59+28=', '81-62='),\n    @('74-32=', '70-24='),\n    @('55+15=', '48-46='),\n    @('5+42=', '25+47='),\n    @('52-24=', '50-11='),\n    @('92-54=', '26+12='),\n    @('84-40=', '93-46='),\n    @('4+0=', '83-41='),\n    @('64+10=', '67-2='),\n    @('34-26=', '31-2='),\n    @('80-13=', '23+1='),\n    @('87-5=', '37+30='),\n    @('74-70=', '67+32='),\n    @('38-22=', '12+24='),\n    @('10+66=', '77-50='),\n    @('71+17=', '15+55='),\n    @('37-31=', '38+1='),\n    @('91-3=', '9+26='),\n    @('28-10=', '82-24='),\n    @('10+40=', '90-20='),\n    @('15-12=', '71-9='),\n    @('96-73=', '85-43='),\n    @('19-1=', '95-87='),\n    @('50-9=', '59-42='),\n    @('26-8=', '28+45='),\n    @('52+42=', '73+3='),\n    @('75-68=', '98-67='),\n    @('75+3=', '70-28='),\n    @('36-15=', '98-73=')\n)\n\nif ($rows * $cols -ne $cellValues.Count) {\n    throw \"Expected $($cellValues.Count) cells, table has $($rows * $cols)\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        $pair = $cellValues[$idx]\n        $oldVal = $pair[0]\n        $newVal = $pair[1]\n        if ($cellText -ne $oldVal) {\n            throw \"Cell ($r,$c) text mismatch: expected [$oldVal], got [$cellText]\"\n        }\n        $cell.Range.Text = $newVal\n        $idx++\n    }\n}\n\n"}
